$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.252.58"
$ws.Range("E2").Value = "  -5.59%  "

$ws.Range("D3").Value = "1.837.79"
$ws.Range("E3").Value = "  -5.23%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.003"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.30%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "330.44"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.54%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.002"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.35%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4596"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -5.02%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3858"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -6.53%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "45.97"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.08%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07861"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.87%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.9657"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -5.00%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "21.90"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -7.42%  "

$ws.Range("D13").Value = "1.865.29"
$ws.Range("E13").Value = "  -5.06%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.713"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -6.26%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.913"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -5.34%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.06870"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.21%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.003"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.44%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "86.78"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -4.84%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000009959"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.95%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "16.91"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -5.03%  "

$ws.Range("E21").Value = "  -0.51%  "

$ws.Range("D22").Value = "28.294.52"
$ws.Range("E22").Value = "  -5.42%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.337"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -5.36%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "10.97"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -7.74%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.149"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.67%  "

$ws.Range("D26").Value = "2.078.11"
$ws.Range("E26").Value = "  -4.91%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "153.86"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.76%  "

$ws.Range("E28").Value = "  -4.32%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.772"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -13.85%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.981"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -5.84%  "

$ws.Range("E31").Value = "  -3.70%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.9414"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -6.78%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09308"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.39%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.285"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -5.52%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.446"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.79%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.324"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -7.12%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06018"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -8.39%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02152"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -5.94%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.150"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -5.08%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.001"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.46%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "7.600"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -5.04%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.5620"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -5.91%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "10.01"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -6.93%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.1783"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.73%  "

$ws.Range("E45").Value = "  -2.76%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.274"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -8.94%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "11.64"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -6.12%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.5291"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -5.09%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.07024"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -6.32%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.838"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -7.71%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "112.81"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.85%  "
